{"js": "// Overview sentence: \"A new user enters their username, email, password and\n// the data is validated and stored, then the user is given a confirmation\n// about their registration.\"\n// becomes: \"A new user enters their username, password and fullname the\n// data is validated and stored, then the user is given a confirmation\n// about their registration.\"\n//\n// i.e. remove \"email, \" and insert \" fullname\" right after \"password and\".\n\nconst body = context.document.body;\n\n// 1) Remove \"email, \" (only occurrence in the document).\nlet results = body.search(\"email, \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Insert \" fullname\" right after \"password and\" (only occurrence).\nresults = body.search(\"password and\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace the matched text with itself plus the new word so the\n  // surrounding run/formatting (sz 24 / szCs 24) is preserved.\n  results.items[0].insertText(\"password and fullname\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Overview sentence: \"A new user enters their username, email, password and\n# the data is validated and stored, then the user is given a confirmation\n# about their registration.\"\n# becomes: \"A new user enters their username, password and fullname the\n# data is validated and stored, then the user is given a confirmation\n# about their registration.\"\n#\n# i.e. remove \"email, \" and insert \" fullname\" right after \"password and\".\n\n$d = $word.ActiveDocument\n\n# 1) Remove \"email, \" (only occurrence in the document).\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"email, \"\n$find1.Replacement.ClearFormatting()\n$find1.Replacement.Text = \"\"\n$find1.Execute(\n  [ref]$find1.Text,          # FindText\n  $false,                    # MatchCase\n  $false,                    # MatchWholeWord\n  $false,                    # MatchWildcards\n  $false,                    # MatchSoundsLike\n  $false,                    # MatchAllWordForms\n  $true,                     # Forward\n  1,                         # Wrap (wdFindContinue)\n  $false,                    # Format\n  $find1.Replacement.Text,   # ReplaceWith\n  2                          # Replace (wdReplaceAll)\n) | Out-Null\n\n# 2) Insert \" fullname\" right after \"password and\" (only occurrence).\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"password and\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"password and fullname\"\n$find2.Execute(\n  [ref]$find2.Text,\n  $false,\n  $false,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  $find2.Replacement.Text,\n  2\n) | Out-Null\n"}
